$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.561.58'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.961.05'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.68'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.378'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0809'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.07%  '
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.34%  '
$ws.Range("D13").Value = '2.248.71'
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.826'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = '1.965.02'
$ws.Range("E17").Value = '  +1.54%  '
$ws.Range("D18").Value = '36.455.89'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").Value = '0.0₃0857'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("E27").Value = '  +3.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0618'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +5.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.29%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0983'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.02'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").Value = '1.363.16'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").Value = '2.140.10'
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.61'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.94%  '
